# The sheet currently has an explicit <cols> block (custom column widths for
# columns A-E). The target revision has no column-width overrides at all, so
# rather than merely resetting ColumnWidth (which still leaves an explicit,
# serialized <col customWidth="1"/> entry behind), we rebuild the sheet from
# a brand-new worksheet that has never had its column formatting touched.
#
# Plan:
#   1. Remember the original sheet's name and position.
#   2. Add a fresh worksheet immediately before it (fresh sheets carry no
#      column overrides).
#   3. Copy every used cell's value over to the fresh sheet.
#   4. Delete the original sheet and rename the fresh one back to the
#      original name, so the workbook ends up with a single sheet in the
#      same place, but with default column formatting.
#   5. Apply the actual content edit: A1 "key" -> "keys".

$wb = $excel.ActiveWorkbook
$oldName = $wb.ActiveSheet.Name

$oldRef = $wb.Worksheets.Item($oldName)
$newRef = $wb.Worksheets.Add($oldRef)
$newName = $newRef.Name

# Re-fetch by name: this runtime resolves sheet handles by position, and
# inserting a sheet shifts indices, so stale handles can silently point at
# the wrong sheet afterwards.
$old = $wb.Worksheets.Item($oldName)
$new = $wb.Worksheets.Item($newName)

$usedRows = $old.UsedRange.Rows.Count
$usedCols = $old.UsedRange.Columns.Count
for ($r = 1; $r -le $usedRows; $r++) {
  for ($c = 1; $c -le $usedCols; $c++) {
    $srcCell = $old.Cells.Item($r, $c)
    $dstCell = $new.Cells.Item($r, $c)
    $dstCell.Value = $srcCell.Value()
  }
}

$toDelete = $wb.Worksheets.Item($oldName)
$toDelete.Delete() | Out-Null

$finalRef = $wb.Worksheets.Item($newName)
$finalRef.Name = $oldName

# Apply the content change described by the diff: A1 "key" -> "keys".
$ws = $wb.Worksheets.Item($oldName)
$ws.Range("A1").Value = "keys"
